$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = [double]"22.24000000000004"
$ws.Range("H2").Value = [double]"9.154334437033729e-05"
$ws.Range("I2").Value = [double]"9.154334437033729e-05"
$ws.Range("L2").Value = [double]"48.53150760069907"
$ws.Range("M2").Value = "[25.83667600299114, 71.226339198407]"
$ws.Range("N2").Value = [double]"8.856105804122549e-05"
$ws.Range("O2").Value = [double]"8.856105804122549e-05"
$ws.Range("P2").Value = [double]"1.540921321580579"
$ws.Range("Q2").Value = "[0.9622896416401163, 2.1195530015210418]"
$ws.Range("R2").Value = [double]"2.717095247461998e-06"
$ws.Range("S2").Value = [double]"2.717095247461998e-06"
$ws.Range("T2").Value = [double]"56.27130688482907"
$ws.Range("U2").Value = "[42.299605481723106, 70.24300828793504]"
$ws.Range("V2").Value = [double]"2.370961205144795e-10"
$ws.Range("W2").Value = [double]"2.370961205144795e-10"
$ws.Range("X2").Value = [double]"16.78574574574577"
$ws.Range("Y2").Value = [double]"14.73761761761764"
$ws.Range("Z2").Value = [double]"18.83387387387391"
$ws.Range("F3").Value = [double]"22.24000000000004"
$ws.Range("H3").Value = [double]"2.553331110755064e-05"
$ws.Range("I3").Value = [double]"2.553331110755064e-05"
$ws.Range("L3").Value = [double]"54.01519777872625"
$ws.Range("M3").Value = "[26.48698105159592, 81.54341450585657]"
$ws.Range("N3").Value = [double]"0.0002707315338124072"
$ws.Range("O3").Value = [double]"0.0002707315338124072"
$ws.Range("P3").Value = [double]"2.157289850212811"
$ws.Range("Q3").Value = "[1.616395018964118, 2.6981846814615036]"
$ws.Range("R3").Value = [double]"3.086422228903984e-10"
$ws.Range("S3").Value = [double]"3.086422228903984e-10"
$ws.Range("T3").Value = [double]"60.8628913922005"
$ws.Range("U3").Value = "[46.04597757105012, 75.67980521335087]"
$ws.Range("V3").Value = [double]"1.38424161022499e-10"
$ws.Range("W3").Value = [double]"1.38424161022499e-10"
$ws.Range("X3").Value = [double]"14.60404404404407"
$ws.Range("Y3").Value = [double]"12.68948948948951"
$ws.Range("Z3").Value = [double]"16.51859859859863"
$ws.Range("F4").Value = [double]"22.24000000000004"
$ws.Range("H4").Value = [double]"0.02407196989342808"
$ws.Range("I4").Value = [double]"0.02407196989342808"
$ws.Range("L4").Value = [double]"34.08051406288463"
$ws.Range("M4").Value = "[2.5626662464085825, 65.59836187936068]"
$ws.Range("N4").Value = [double]"0.03469717000190253"
$ws.Range("O4").Value = [double]"0.03469717000190253"
$ws.Range("P4").Value = [double]"2.144710900648888"
$ws.Range("Q4").Value = "[0.8616580451287321, 3.4277637561690435]"
$ws.Range("R4").Value = [double]"0.00156572883028061"
$ws.Range("S4").Value = [double]"0.00156572883028061"
$ws.Range("T4").Value = [double]"58.2063749752843"
$ws.Range("U4").Value = "[41.16959939514663, 75.24315055542196]"
$ws.Range("V4").Value = [double]"1.529049686865847e-08"
$ws.Range("W4").Value = [double]"1.529049686865847e-08"
$ws.Range("X4").Value = [double]"14.64856856856859"
$ws.Range("Y4").Value = [double]"10.10706706706708"
$ws.Range("Z4").Value = [double]"19.1900700700701"
$ws.Range("F5").Value = [double]"22.24000000000004"
$ws.Range("H5").Value = [double]"7.474730144763075e-05"
$ws.Range("I5").Value = [double]"7.474730144763075e-05"
$ws.Range("L5").Value = [double]"48.41606867681505"
$ws.Range("M5").Value = "[24.871077760180412, 71.96105959344969]"
$ws.Range("N5").Value = [double]"0.0001496920127117551"
$ws.Range("O5").Value = [double]"0.0001496920127117551"
$ws.Range("P5").Value = [double]"-3.069263693597235"
$ws.Range("Q5").Value = "[-3.6856322222294673, -2.4528951649650037]"
$ws.Range("R5").Value = [double]"4.742872761198669e-13"
$ws.Range("S5").Value = [double]"4.742872761198669e-13"
$ws.Range("T5").Value = [double]"54.92644919653748"
$ws.Range("U5").Value = "[40.320303876261946, 69.53259451681302]"
$ws.Range("V5").Value = [double]"1.446968767027101e-09"
$ws.Range("W5").Value = [double]"1.446968767027101e-09"
$ws.Range("X5").Value = [double]"10.863983983984"
$ws.Range("Y5").Value = [double]"8.682282282282298"
$ws.Range("Z5").Value = [double]"13.04568568568571"
$ws.Range("F6").Value = [double]"24.78000000000043"
$ws.Range("H6").Value = [double]"8.796405259747786e-08"
$ws.Range("I6").Value = [double]"8.796405259747786e-08"
$ws.Range("L6").Value = [double]"73.73081111289808"
$ws.Range("M6").Value = "[49.68005702762264, 97.78156519817352]"
$ws.Range("N6").Value = [double]"1.718867539413083e-07"
$ws.Range("O6").Value = [double]"1.718867539413083e-07"
$ws.Range("P6").Value = [double]"-2.981211046649774"
$ws.Range("Q6").Value = "[-3.3585795335674664, -2.6038425597320813]"
$ws.Range("R6").Value = 0
$ws.Range("S6").Value = 0
$ws.Range("T6").Value = [double]"70.54352598952825"
$ws.Range("U6").Value = "[55.589162509969945, 85.49788946908654]"
$ws.Range("V6").Value = [double]"2.520650355108955e-12"
$ws.Range("W6").Value = [double]"2.520650355108955e-12"
$ws.Range("X6").Value = [double]"11.75747747747768"
$ws.Range("Y6").Value = [double]"10.26918918918937"
$ws.Range("Z6").Value = [double]"13.245765765766"
$ws.Range("F7").Value = [double]"24.78000000000043"
$ws.Range("H7").Value = [double]"0.002998563651989472"
$ws.Range("I7").Value = [double]"0.002998563651989472"
$ws.Range("L7").Value = [double]"49.26014777220632"
$ws.Range("M7").Value = "[17.075661257903334, 81.44463428650931]"
$ws.Range("N7").Value = [double]"0.003497267096751022"
$ws.Range("O7").Value = [double]"0.003497267096751022"
$ws.Range("P7").Value = [double]"-2.993789996213697"
$ws.Range("Q7").Value = "[-3.849158566560468, -2.138421425866926]"
$ws.Range("R7").Value = [double]"8.609838841877604e-09"
$ws.Range("S7").Value = [double]"8.609838841877604e-09"
$ws.Range("T7").Value = [double]"73.56746559175026"
$ws.Range("U7").Value = "[54.63169299131479, 92.50323819218573]"
$ws.Range("V7").Value = [double]"6.203901836698833e-10"
$ws.Range("W7").Value = [double]"6.203901836698833e-10"
$ws.Range("X7").Value = [double]"11.80708708708729"
$ws.Range("Y7").Value = [double]"8.433633633633779"
$ws.Range("Z7").Value = [double]"15.18054054054081"
$ws.Range("F8").Value = [double]"24.78000000000043"
$ws.Range("H8").Value = [double]"7.168182212358154e-05"
$ws.Range("I8").Value = [double]"7.168182212358154e-05"
$ws.Range("L8").Value = [double]"57.19566578990999"
$ws.Range("M8").Value = "[28.020218137313762, 86.37111344250621]"
$ws.Range("N8").Value = [double]"0.0002737395459013747"
$ws.Range("O8").Value = [double]"0.0002737395459013747"
$ws.Range("P8").Value = [double]"2.912026824048196"
$ws.Range("Q8").Value = "[2.3082372449798885, 3.5158164031165042]"
$ws.Range("R8").Value = [double]"1.280975325812506e-12"
$ws.Range("S8").Value = [double]"1.280975325812506e-12"
$ws.Range("T8").Value = [double]"63.5949403875056"
$ws.Range("U8").Value = "[47.10161801074365, 80.08826276426755]"
$ws.Range("V8").Value = [double]"7.568081716868846e-10"
$ws.Range("W8").Value = [double]"7.568081716868846e-10"
$ws.Range("X8").Value = [double]"13.29537537537561"
$ws.Range("Y8").Value = [double]"10.91411411411431"
$ws.Range("Z8").Value = [double]"15.67663663663691"
$ws.Range("F9").Value = [double]"24.78000000000043"
$ws.Range("H9").Value = [double]"0.0005108293604034753"
$ws.Range("I9").Value = [double]"0.0005108293604034753"
$ws.Range("L9").Value = [double]"48.58625898128606"
$ws.Range("M9").Value = "[16.786409214089787, 80.38610874848233]"
$ws.Range("N9").Value = [double]"0.003549843497206329"
$ws.Range("O9").Value = [double]"0.003549843497206329"
$ws.Range("P9").Value = [double]"2.610132034514042"
$ws.Range("Q9").Value = "[1.9811845563178876, 3.2390795127101963]"
$ws.Range("R9").Value = [double]"1.042734787404243e-10"
$ws.Range("S9").Value = [double]"1.042734787404243e-10"
$ws.Range("T9").Value = [double]"47.97950310899763"
$ws.Range("U9").Value = "[31.731876156904875, 64.22713006109039]"
$ws.Range("V9").Value = [double]"3.733610369760498e-07"
$ws.Range("W9").Value = [double]"3.733610369760498e-07"
$ws.Range("X9").Value = [double]"14.48600600600626"
$ws.Range("Y9").Value = [double]"12.00552552552574"
$ws.Range("Z9").Value = [double]"16.96648648648678"
$ws.Range("F10").Value = [double]"24.78000000000043"
$ws.Range("H10").Value = [double]"8.901693340157557e-05"
$ws.Range("I10").Value = [double]"8.901693340157557e-05"
$ws.Range("L10").Value = [double]"51.59817463921517"
$ws.Range("M10").Value = "[22.191825662576917, 81.00452361585343]"
$ws.Range("N10").Value = [double]"0.0009600255930237722"
$ws.Range("O10").Value = [double]"0.0009600255930237722"
$ws.Range("P10").Value = [double]"2.081816152829273"
$ws.Range("Q10").Value = "[1.490605523324887, 2.6730267823336584]"
$ws.Range("R10").Value = [double]"7.438825111449887e-09"
$ws.Range("S10").Value = [double]"7.438825111449887e-09"
$ws.Range("T10").Value = [double]"57.34279528399195"
$ws.Range("U10").Value = "[41.856767078387975, 72.82882348959592]"
$ws.Range("V10").Value = [double]"2.143939248355764e-09"
$ws.Range("W10").Value = [double]"2.143939248355764e-09"
$ws.Range("X10").Value = [double]"16.5696096096099"
$ws.Range("Y10").Value = [double]"14.2379579579582"
$ws.Range("Z10").Value = [double]"18.90126126126159"
$ws.Range("F11").Value = [double]"24.78000000000043"
$ws.Range("H11").Value = [double]"0.01464160680763471"
$ws.Range("I11").Value = [double]"0.01464160680763471"
$ws.Range("L11").Value = [double]"34.20458127284717"
$ws.Range("M11").Value = "[3.3406095942790586, 65.06855295141528]"
$ws.Range("N11").Value = [double]"0.03063096514141916"
$ws.Range("O11").Value = [double]"0.03063096514141916"
$ws.Range("P11").Value = [double]"2.446605690183042"
$ws.Range("Q11").Value = "[1.6541318676558872, 3.2390795127101972]"
$ws.Range("R11").Value = [double]"1.480331108805188e-07"
$ws.Range("S11").Value = [double]"1.480331108805188e-07"
$ws.Range("T11").Value = [double]"52.23666730399489"
$ws.Range("U11").Value = "[36.21563632257832, 68.25769828541146]"
$ws.Range("V11").Value = [double]"4.481566495684319e-08"
$ws.Range("W11").Value = [double]"4.481566495684319e-08"
$ws.Range("X11").Value = [double]"15.1309309309312"
$ws.Range("Y11").Value = [double]"12.00552552552574"
$ws.Range("Z11").Value = [double]"18.25633633633666"
$ws.Range("F12").Value = [double]"24.78000000000043"
$ws.Range("H12").Value = [double]"0.004914660404543469"
$ws.Range("I12").Value = [double]"0.004914660404543469"
$ws.Range("L12").Value = [double]"46.08514379726368"
$ws.Range("M12").Value = "[10.568541372691854, 81.6017462218355]"
$ws.Range("N12").Value = [double]"0.01214622087371309"
$ws.Range("O12").Value = [double]"0.01214622087371309"
$ws.Range("P12").Value = [double]"2.55981623625835"
$ws.Range("Q12").Value = "[1.7925003128590413, 3.3271321596576593]"
$ws.Range("R12").Value = [double]"2.661531928183081e-08"
$ws.Range("S12").Value = [double]"2.661531928183081e-08"
$ws.Range("T12").Value = [double]"71.02990520480036"
$ws.Range("U12").Value = "[52.106844200727494, 89.95296620887322]"
$ws.Range("V12").Value = [double]"1.516504477550029e-09"
$ws.Range("W12").Value = [double]"1.516504477550029e-09"
$ws.Range("X12").Value = [double]"14.6844444444447"
$ws.Range("Y12").Value = [double]"11.65825825825846"
$ws.Range("Z12").Value = [double]"17.71063063063094"
$ws.Range("F13").Value = [double]"24.78000000000043"
$ws.Range("H13").Value = [double]"0.0004061518540575682"
$ws.Range("I13").Value = [double]"0.0004061518540575682"
$ws.Range("L13").Value = [double]"50.10716641257162"
$ws.Range("M13").Value = "[19.924325492639497, 80.29000733250373]"
$ws.Range("N13").Value = [double]"0.001673418879329924"
$ws.Range("O13").Value = [double]"0.001673418879329924"
$ws.Range("P13").Value = [double]"1.968605606753965"
$ws.Range("Q13").Value = "[1.3019212798660407, 2.6352899336418885]"
$ws.Range("R13").Value = [double]"3.738387701623935e-07"
$ws.Range("S13").Value = [double]"3.738387701623935e-07"
$ws.Range("T13").Value = [double]"71.4493238388888"
$ws.Range("U13").Value = "[54.68565586550149, 88.21299181227612]"
$ws.Range("V13").Value = [double]"4.939959552530127e-11"
$ws.Range("W13").Value = [double]"4.939959552530127e-11"
$ws.Range("X13").Value = [double]"17.01609609609639"
$ws.Range("Y13").Value = [double]"14.38678678678704"
$ws.Range("Z13").Value = [double]"19.64540540540575"
$ws.Range("B14").Value = 0
$ws.Range("F14").Value = [double]"24.78000000000043"
$ws.Range("H14").Value = [double]"0.06203153348447354"
$ws.Range("I14").Value = [double]"0.06203153348447354"
$ws.Range("L14").Value = [double]"31.09141693338714"
$ws.Range("M14").Value = "[-1.2333142909686643, 63.41614815774294]"
$ws.Range("N14").Value = [double]"0.05900616465077935"
$ws.Range("O14").Value = [double]"0.05900616465077935"
$ws.Range("P14").Value = [double]"1.440289725069194"
$ws.Range("Q14").Value = "[-0.01886842434588587, 2.8994478744842738]"
$ws.Range("R14").Value = [double]"0.05290919250862869"
$ws.Range("S14").Value = [double]"0.05290919250862869"
$ws.Range("T14").Value = [double]"56.67034556902512"
$ws.Range("U14").Value = "[38.210079519481226, 75.130611618569]"
$ws.Range("V14").Value = [double]"1.66956434455301e-07"
$ws.Range("W14").Value = [double]"1.66956434455301e-07"
$ws.Range("X14").Value = [double]"19.09969969970004"
$ws.Range("Y14").Value = [double]"13.34498498498522"
$ws.Range("Z14").Value = [double]"24.85441441441485"

Write-Host "Applied 235 cell updates"